$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change C2 from text "hard" to numeric value 0
$ws.Range("C2").Value = 0

# Update the selected cell/range to C8
$ws.Range("C8").Select()
